# Add "Area" column (G/H) and a condensed Qtotal/Atotal summary (J/K) to
# Sheet1, mirroring the existing "Q" discharge calculation with a parallel
# cross-sectional-area calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: first area segment + running totals ---------------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Condensed summary pulling the two totals together
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Row 3: second area segment (standalone formula) -----------------------
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- Rows 4-15: remaining area segments as one shared formula -------------
# (row 15 is new - the area series runs one row past the existing data)
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Match the selection left behind by the edit
$ws.Range("J2:K2").Select() | Out-Null
